$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("2024-05-17 16:30:04","NVDA","BUY","STRONG_BUY","BUY"),
    @("2024-05-17 16:30:04","UNH","BUY","STRONG_BUY","SELL"),
    @("2024-05-17 16:30:04","GOOG","SELL","STRONG_SELL","NEUTRAL"),
    @("2024-05-17 16:30:04","SPY","BUY","BUY","NEUTRAL"),
    @("2024-05-17 16:30:04","GOLD","BUY","STRONG_BUY","NEUTRAL"),
    @("2024-05-17 16:30:04","AMZN","NEUTRAL","SELL","BUY"),
    @("2024-05-17 16:30:04","BABA","STRONG_BUY","STRONG_BUY","BUY"),
    @("2024-05-17 16:30:04","META","SELL","STRONG_SELL","BUY"),
    @("2024-05-17 16:30:04","DOCN","SELL","STRONG_SELL","NEUTRAL"),
    @("2024-05-17 16:30:04","VOO","NEUTRAL","BUY","NEUTRAL"),
    @("2024-05-17 16:30:04","TSLA","BUY","STRONG_BUY","NEUTRAL"),
    @("2024-05-17 16:30:04","MSFT","SELL","SELL","NEUTRAL"),
    @("2024-05-17 16:30:04","AAPL","BUY","STRONG_BUY","SELL"),
    @("2024-05-17 16:45:04","NVDA","BUY","STRONG_BUY","BUY"),
    @("2024-05-17 16:45:04","UNH","BUY","STRONG_BUY","SELL"),
    @("2024-05-17 16:45:04","GOOG","SELL","STRONG_SELL","NEUTRAL"),
    @("2024-05-17 16:45:04","SPY","BUY","BUY","NEUTRAL"),
    @("2024-05-17 16:45:04","GOLD","BUY","STRONG_BUY","NEUTRAL"),
    @("2024-05-17 16:45:04","AMZN","NEUTRAL","SELL","BUY"),
    @("2024-05-17 16:45:04","BABA","STRONG_BUY","STRONG_BUY","BUY"),
    @("2024-05-17 16:45:04","META","SELL","STRONG_SELL","BUY"),
    @("2024-05-17 16:45:04","DOCN","SELL","STRONG_SELL","NEUTRAL"),
    @("2024-05-17 16:45:04","VOO","NEUTRAL","BUY","NEUTRAL"),
    @("2024-05-17 16:45:04","TSLA","BUY","STRONG_BUY","NEUTRAL"),
    @("2024-05-17 16:45:04","MSFT","SELL","SELL","NEUTRAL"),
    @("2024-05-17 16:45:04","AAPL","BUY","STRONG_BUY","SELL"),
    @("2024-05-17 17:00:04","NVDA","BUY","STRONG_BUY","BUY"),
    @("2024-05-17 17:00:04","UNH","BUY","STRONG_BUY","SELL"),
    @("2024-05-17 17:00:04","GOOG","SELL","STRONG_SELL","NEUTRAL"),
    @("2024-05-17 17:00:04","SPY","BUY","BUY","NEUTRAL"),
    @("2024-05-17 17:00:04","GOLD","BUY","STRONG_BUY","NEUTRAL"),
    @("2024-05-17 17:00:04","AMZN","NEUTRAL","SELL","BUY"),
    @("2024-05-17 17:00:04","BABA","STRONG_BUY","STRONG_BUY","BUY"),
    @("2024-05-17 17:00:04","META","SELL","STRONG_SELL","BUY"),
    @("2024-05-17 17:00:04","DOCN","SELL","STRONG_SELL","NEUTRAL"),
    @("2024-05-17 17:00:04","VOO","NEUTRAL","BUY","NEUTRAL"),
    @("2024-05-17 17:00:04","TSLA","BUY","STRONG_BUY","NEUTRAL"),
    @("2024-05-17 17:00:04","MSFT","SELL","SELL","NEUTRAL"),
    @("2024-05-17 17:00:04","AAPL","BUY","STRONG_BUY","SELL")
)

$startRow = 184
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = ""
}
